# Fix bug in balance_type creation: the District_Heating connection (row 5 on the
# Connections sheet) was missing its mirrored Input2/Output2 values and was left as
# "connection_type_normal" instead of "connection_type_lossless_bidirectional",
# which left the Waste_Heat source node unbalanced and made the model infeasible.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Connections")

# Mirror Input2 (C5) with Output1 (D5), and Output2 (E5) with Input1 (B5),
# matching the pattern used by every other connection row.
$ws.Range("C5").Value = "District_Heating"
$ws.Range("E5").Value = "Waste_Heat"

# Correct the connection type so the District_Heating connection is balanced like
# the rest of the bidirectional connections.
$ws.Range("F5").Value = "connection_type_lossless_bidirectional"

# Reflect the cell that was last active after making the edit.
$ws.Range("G14").Select()
